$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.188.32"
$ws.Range("D3").Value = "1.589.97"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'211.78"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.813.86"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.585.33"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "'63.62"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "26.173.41"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "'214.66"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").Value = "'7.35"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").Value = "'144.83"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").Value = "'0.0493"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "1.418.74"
$ws.Range("E33").Value = "  +8.09%  "
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("E37").Value = "  -3.94%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("E40").Value = "  +4.70%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'0.963"
$ws.Range("E42").Value = "  -11.99%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "1.725.61"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'61.04"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "'87.07"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "'0.0959"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("E51").Value = "  -0.15%  "
